# Edit script: reshape the "Input" sheet to the new standard template
# column layout, add company/delivery e-mail columns, drop the old
# 발주번호/단위/공급가액/부가세 columns, rename headers, and clean up
# the leftover empty "비고" cells on the 갑지/을지 sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Input" - full reshape to the new template column order
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Input")

$headers = @(
    "발주일자", "납기일자", "거래처명", "거래처 이메일", "납품처명",
    "납품처 이메일", "프로젝트명", "대분류", "중분류", "소분류",
    "품목명", "규격", "수량", "단가", "총금액", "비고"
)

$row2 = @(
    "2025-09-15", "2025-09-14", "케이에스파워텍", "케이에스파워텍@example.com",
    "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차",
    "5. 운반비", "일반자재", "기타", "2월 운반비", "KS규격-1",
    1, 0, 0, ""
)

$row3 = @(
    "2025-09-08", "2025-09-12", "케이에스파워텍", "케이에스파워텍@example.com",
    "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차",
    "2. 부자재비", "3) 기타", "기타", "스텐망 663*1670", "KS규격-2",
    2, 29000, 63800, ""
)

$row4 = @(
    "2025-09-03", "2025-10-09", "케이에스파워텍", "케이에스파워텍@example.com",
    "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차",
    "2. 부자재비", "3) 기타", "기타", "스텐망 1023*1100", "KS규격-3",
    2, 29000, 63800, ""
)

# Clear the whole previous footprint (A1:Q4) before laying out the new grid.
# Formats are cleared first (this also removes the bold/bordered header
# style from the old template) so that later ClearContents calls on
# still-untouched cells do not leave behind stray empty <c> elements.
$ws1.Range("A1:Q4").ClearContents()
$ws1.Range("A1:Q4").ClearFormats()

# The date-looking text (발주일자/납기일자) must stay plain text, not get
# auto-converted into Excel date serials - force a text number format
# on those two columns before writing the values.
$ws1.Range("A2:B4").NumberFormat = "@"

for ($col = 1; $col -le $headers.Length; $col++) {
    $ws1.Cells.Item(1, $col).Value = $headers[$col - 1]
}
for ($col = 1; $col -le $row2.Length; $col++) {
    $ws1.Cells.Item(2, $col).Value = $row2[$col - 1]
}
for ($col = 1; $col -le $row3.Length; $col++) {
    $ws1.Cells.Item(3, $col).Value = $row3[$col - 1]
}
for ($col = 1; $col -le $row4.Length; $col++) {
    $ws1.Cells.Item(4, $col).Value = $row4[$col - 1]
}

# Now that the text values are safely stored as strings, drop the
# temporary text number format again so the cells end up unstyled,
# matching the rest of the (unstyled) data rows.
$ws1.Range("A2:B4").ClearFormats()

# The old blank "비고" values should stay truly empty cells (no content),
# matching the rest of the workbook's treatment of blank trailing columns.
$ws1.Range("P2:P4").ClearContents()

# Drop the now-unused column Q entirely (old last column, beyond new P).
$ws1.Range("Q1:Q4").ClearContents()

# ---------------------------------------------------------------
# Sheets 2 & 3: "갑지" / "을지" - blank out the leftover empty 비고 cells
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("갑지")
$ws2.Range("I2:I4").ClearContents()

$ws3 = $wb.Worksheets.Item("을지")
$ws3.Range("I2:I4").ClearContents()
